$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.214.96'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '1.786.27'
$ws.Range("E3").Value = '  -0.31%  '

$ws.Range("E4").Value = '  +0.18%  '

$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'225.72"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("E9").Value = '  +0.02%  '

$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.0687"
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  -0.79%  '

$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("D12").Value = '2.044.54'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.803.49'
$ws.Range("E13").Value = '  +0.64%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'10.96"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = '  -5.02%  '

$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.626"
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = '  +0.46%  '

$ws.Range("D16").Value = '34.174.34'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("E17").Value = '  -0.38%  '

$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'67.93"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("E19").Value = '  +2.48%  '

$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").Value = "'246.04"
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("E23").Value = '  +1.07%  '

$ws.Range("E24").Value = '  +0.23%  '

$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'161.76"
$ws.Range("D25").Style = $style_D25

$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("E28").Value = '  +1.34%  '

$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("E30").Value = '  -0.99%  '

$ws.Range("E31").Value = '  -0.05%  '

$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.75"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("E33").Value = '  +3.78%  '

$ws.Range("E34").Value = '  -1.57%  '

$ws.Range("D35").Value = '1.439.85'
$ws.Range("E35").Value = '  -0.90%  '

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'2.61"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  +10.44%  '

$ws.Range("E37").Value = '  +2.21%  '

$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("E39").Value = '  -1.38%  '

$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'81.60"
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = '  +1.32%  '

$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'2.40"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  +1.95%  '

$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'14.07"
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = '  +5.36%  '

$ws.Range("E43").Value = '  +1.25%  '

$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("E45").Value = '  +2.12%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").Value = "'6.05"
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = '  -0.39%  '

$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'1.08"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  +0.86%  '

$ws.Range("D48").Value = '1.940.60'
$ws.Range("E48").Value = '  -0.49%  '

$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'105.41"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = '  -2.15%  '

$ws.Range("E50").Value = '  +0.15%  '
